# Weekly driver report update for 2025-04-19
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text string into a cell without Excel's
# autodetection turning date-shaped text (e.g. "2024-11-10") into a real
# date (which would also mint a new number-format style on the cell).
# We stage the text - forced to Text format - in a scratch cell that is
# already inside the sheet's used range, then paste-special *values only*
# into the destination so the destination keeps its original style.
function Set-TextValue {
    param($destRange, [string]$text)

    $scratch = $ws.Range("A20")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $destRange.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# --- Bad Drivers table ---
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.90.0.5"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 153
# D3 (Good Roaming Calculation %) is unchanged

# --- Totals row ---
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 153

# --- Good Drivers table (rows 12-17) ---
# Row 12
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B12").Value = 56018
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = 0

# Row 13
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B13").Value = 34244
$ws.Range("D13").Value = 100
$ws.Range("E13").Value = 0

# Row 14
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B14").Value = 442178
$ws.Range("D14").Value = 99.90000000000001
Set-TextValue $ws.Range("E14") "2024-11-10"

# Row 15
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B15").Value = 77849
$ws.Range("D15").Value = 99.90000000000001
Set-TextValue $ws.Range("E15") "2021-08-18"

# Row 16
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B16").Value = 59673
$ws.Range("D16").Value = 100
Set-TextValue $ws.Range("E16") "2020-08-05"

# Row 17
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B17").Value = 113652
$ws.Range("D17").Value = 100
Set-TextValue $ws.Range("E17") "2019-12-14"
